$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force new cells to be stored as text (matches the source data's inlineStr
# string typing, even for numeric-looking values like SKUs/prices) without
# leaving a lingering custom number-format style behind.
$textRange = $ws.Range("A8:E13")
$textRange.NumberFormat = "@"

$ws.Range("A8").Value = "764804"
$ws.Range("B8").Value = "48 SparClean Dish Manual"
$ws.Range("C8").Value = "1"
$ws.Range("D8").Value = "102.51"
$ws.Range("E8").Value = "102.51"

$ws.Range("A9").Value = "T607646"
$ws.Range("B9").Value = "Container - Deli (64oz)"
$ws.Range("C9").Value = "1"
$ws.Range("D9").Value = "123.23"
$ws.Range("E9").Value = "123.23"

# A10 has no SKU in the source data - represented as an explicit (empty)
# text cell rather than a truly blank cell.
$ws.Range("A10").Value = "'"

$ws.Range("B10").Value = "Lid - Deli (64oz)"
$ws.Range("C10").Value = "1"
$ws.Range("D10").Value = "82.63"
$ws.Range("E10").Value = "82.63"

$ws.Range("A11").Value = "LKC1220F"
$ws.Range("B11").Value = "Lid Cold Flat - 12/20oz (No Slot)"
$ws.Range("C11").Value = "1"
$ws.Range("D11").Value = "48.10"
$ws.Range("E11").Value = "48.10"

$ws.Range("A12").Value = "711603"
$ws.Range("B12").Value = "NABC Bathroom Cleaner"
$ws.Range("C12").Value = "1"
$ws.Range("D12").Value = "35.35"
$ws.Range("E12").Value = "35.35"

$ws.Range("A13").Value = "LKC1624F"
$ws.Range("B13").Value = "Parfait Lid"
$ws.Range("C13").Value = "1"
$ws.Range("D13").Value = "47.53"
$ws.Range("E13").Value = "47.53"

# Restore default styling so the newly written cells don't carry a custom
# number-format style index (source cells had no explicit style either).
$textRange.Style = "Normal"
